# Fix bugs in "print_to_excel" and "calculate_time_cost_per_group":
# Update the "address" column (D) values on the "Child" sheet for rows 2-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "-3,2"
    3  = "-2,1"
    4  = "-10,-9"
    5  = "-3,0"
    6  = "-4,-1"
    7  = "-10,1"
    8  = "-1,7"
    9  = "2,4"
    10 = "-6,-9"
    11 = "-4,-3"
    12 = "-7,1"
    13 = "-7,-5"
    14 = "-3,0"
    15 = "-3,7"
    16 = "-5,-1"
    17 = "-8,-1"
    18 = "-1,-8"
    19 = "7,8"
    20 = "-3,-6"
    21 = "-8,8"
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
